$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "Tenant_ID" column header to "Tenant Passport ID Number"
$ws.Range("C1").Value = "Tenant Passport ID Number"

# Update the active selection to C2 (as left after editing the header cell)
$ws.Range("C2").Select()
